$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D/E values that look numeric get a leading apostrophe so Excel
# stores them as literal text (matching the original inlineStr cells)
# instead of auto-converting them to numbers/percentages.

$ws.Range("D2").Value = "'309.98"
$ws.Range("E2").Value = "'-0.52%"
$ws.Range("D3").Value = "'36.99"
$ws.Range("E3").Value = "'-1.99%"
$ws.Range("D4").Value = "'5.124"
$ws.Range("E4").Value = "'0.11%"
$ws.Range("D5").Value = "'0.07871"
$ws.Range("E5").Value = "'-0.50%"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").Value = "'8.265"
$ws.Range("E6").Value = "'0.34%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.880"
$ws.Range("E7").Value = "'-1.40%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.943"
$ws.Range("E8").Value = "'-4.42%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9240"
$ws.Range("E9").Value = "'-0.34%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1173"
$ws.Range("E10").Value = "'-2.48%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1900"
$ws.Range("E11").Value = "'-0.41%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08959"
$ws.Range("E12").Value = "'-3.91%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03318"
$ws.Range("E13").Value = "'-1.90%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09606"
$ws.Range("E14").Value = "'-0.08%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001386"
$ws.Range("E15").Value = "'0.97%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006191"
$ws.Range("E16").Value = "'5.52%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.393"
$ws.Range("E17").Value = "'-3.88%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.399"
$ws.Range("E18").Value = "'-0.33%"
$ws.Range("D19").Value = "'0.3459"
$ws.Range("E19").Value = "'0.28%"
$ws.Range("D20").Value = "'6.365"
$ws.Range("E20").Value = "'20.69%"
$ws.Range("D21").Value = "'0.1314"
$ws.Range("E21").Value = "'2.54%"
$ws.Range("D22").Value = "'0.2403"
$ws.Range("E22").Value = "'-7.08%"
$ws.Range("D23").Value = "'0.04345"
$ws.Range("D24").Value = "'0.001199"
$ws.Range("E24").Value = "'-4.01%"
$ws.Range("D25").Value = "'0.004284"
$ws.Range("E25").Value = "'0.23%"
$ws.Range("D26").Value = "'0.0001399"
$ws.Range("E26").Value = "'7.92%"
$ws.Range("D27").Value = "'0.0002898"
$ws.Range("D39").Value = "'0.02160"
$ws.Range("E39").Value = "'3.50%"
$ws.Range("E40").Value = "'-1.15%"
$ws.Range("D41").Value = "'0.007586"
$ws.Range("E41").Value = "'-0.33%"
$ws.Range("D42").Value = "'0.1355"
$ws.Range("E42").Value = "'0.16%"
$ws.Range("D43").Value = "'0.008527"
$ws.Range("E43").Value = "'-6.36%"
$ws.Range("D44").Value = "'0.002069"
$ws.Range("E44").Value = "'-0.76%"
$ws.Range("D45").Value = "'0.008126"
$ws.Range("E45").Value = "'-5.82%"
$ws.Range("D46").Value = "'0.00006558"
$ws.Range("E46").Value = "'-1.86%"
$ws.Range("E47").Value = "'0.19%"
$ws.Range("D48").Value = "'0.003292"
$ws.Range("E48").Value = "'14.20%"
$ws.Range("D49").Value = "'0.001442"
$ws.Range("E49").Value = "'20.43%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'0.19%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'0.19%"
